# PatientH_df.xlsx update
#
# The sheet is a simple 48-row lookup table (A: label, B/C/D: numeric
# measurements). This edit:
#   1. Removes the row whose label is "pc" (originally row 23), shifting
#      every row below it up by one.
#   2. Appends a brand-new row at the bottom of the table with label
#      "zy_r" and its three measurements.
#   3. Leaves the selection on the last cell that was typed into (D48),
#      mirroring how the sheet was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the "pc" row (row 23) - shifts rows 24:48 up to 23:47 and
# Excel automatically drops "pc" from the shared-string table since it
# becomes unused.
$ws.Rows("23:23").Delete()

# 2) Add the new trailing row (now row 48) with its label + values.
$ws.Range("A48").Value = "zy_r"
$ws.Range("B48").Value = -56.73
$ws.Range("C48").Value = 5.33
$ws.Range("D48").Value = 75.06

# 3) Leave the selection where the user finished typing.
$ws.Range("D48").Select()
